# Edit script implementing:
#  - Convert the costs (capex & fom) to EUR-2025 basis: for every data row,
#    set the 2030 (G), 2040 (H) and 2050 (I) columns equal to the 2025 (F)
#    column value.
#  - Remove rows with zero production from the NUTS1 production catalogues
#    (2018 and 30/50 tables), for the block of rows describing CH / NO / UK
#    (rows 152-177), shifting the remaining rows up.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) ind_process_routes_capex: propagate 2025 cost (col F) into
#    2030/2040/2050 (cols G/H/I) for data rows 2..64
# ---------------------------------------------------------------------
$wsCapex = $wb.Worksheets.Item("ind_process_routes_capex")
$lastRowCapex = $wsCapex.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRowCapex; $r++) {
    $f2025 = $wsCapex.Cells.Item($r, 6).Value()
    $wsCapex.Cells.Item($r, 7).Value = $f2025
    $wsCapex.Cells.Item($r, 8).Value = $f2025
    $wsCapex.Cells.Item($r, 9).Value = $f2025
}

# ---------------------------------------------------------------------
# 2) ind_process_routes_fom: same treatment for all data rows
# ---------------------------------------------------------------------
$wsFom = $wb.Worksheets.Item("ind_process_routes_fom")
$lastRowFom = $wsFom.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRowFom; $r++) {
    $f2025 = $wsFom.Cells.Item($r, 6).Value()
    $wsFom.Cells.Item($r, 7).Value = $f2025
    $wsFom.Cells.Item($r, 8).Value = $f2025
    $wsFom.Cells.Item($r, 9).Value = $f2025
}

# ---------------------------------------------------------------------
# 3) ind_production_2018_nuts1: remove rows (152..177) whose production
#    value (column F) is zero, shifting remaining rows up. Rows are
#    collected first, then deleted from the bottom up so row numbers
#    for the remaining deletions stay valid.
# ---------------------------------------------------------------------
$wsProd2018 = $wb.Worksheets.Item("ind_production_2018_nuts1")
$zeroRows2018 = New-Object System.Collections.ArrayList
for ($r = 152; $r -le 177; $r++) {
    $prod = $wsProd2018.Cells.Item($r, 6).Value()
    if ($prod -eq 0) {
        [void]$zeroRows2018.Add($r)
    }
}
for ($i = $zeroRows2018.Count - 1; $i -ge 0; $i--) {
    $wsProd2018.Rows.Item($zeroRows2018[$i]).Delete()
}

# ---------------------------------------------------------------------
# 4) ind_production_30_50_nuts1: same row removal (production value is
#    in column E, mirrored in column F).
# ---------------------------------------------------------------------
$wsProd3050 = $wb.Worksheets.Item("ind_production_30_50_nuts1")
$zeroRows3050 = New-Object System.Collections.ArrayList
for ($r = 152; $r -le 177; $r++) {
    $prod = $wsProd3050.Cells.Item($r, 5).Value()
    if ($prod -eq 0) {
        [void]$zeroRows3050.Add($r)
    }
}
for ($i = $zeroRows3050.Count - 1; $i -ge 0; $i--) {
    $wsProd3050.Rows.Item($zeroRows3050[$i]).Delete()
}

Write-Output "edit complete"
